# Make the body of the Abstract paragraph bold (everything after "Abstract - ").
$d = $word.ActiveDocument

# Locate the Abstract paragraph (it starts with "Abstract - ").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Abstract")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $full = $target.Range.Text
    $idx = $full.IndexOf("There have been")
    if ($idx -ge 0) {
        $pStart = $target.Range.Start
        $pEnd = $target.Range.End

        $bodyStart = $pStart + $idx
        # Exclude the trailing paragraph-mark character from the range.
        $bodyEnd = $pEnd - 1

        $body = $d.Range($bodyStart, $bodyEnd)
        $body.Font.Bold = 1
        $body.Font.BoldBi = 1
    }
}
